# Applies the "update from work laptop" edit:
#  - splits a couple of runs and wraps the affected words in <w:proofErr/>
#    markers (grammar/spelling squiggles Word leaves behind once it has
#    checked text it just saw edited)
#  - appends several new paragraphs of notes after "A(iii)"
#  - leaves a _GoBack bookmark at the very end (Word always drops one at
#    the last edited location)

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph 2 ("Clearly x2 has more importance ...") -------------------
# Split off "5." into its own run wrapped in a gramStart/gramEnd proofErr
# pair (the "937.." run is unchanged but now closes the gram range).
$para2 = $d.Paragraphs.Item(2)
$p2xml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Clearly x2 has more importance (its slope is </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>5.</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>937..</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>. The lager X2 is the more likely it will be a +1.</w:t></w:r>' +
  '</w:p>'
$para2.Range.InsertXML($pkgHeader + $p2xml + $pkgFooter)

# --- Paragraph 3 ("X1 has a slope of ...") ---------------------------------
# Split off "doesnt" into its own run wrapped in a spellStart/spellEnd pair.
$para3 = $d.Paragraphs.Item(3)
$p3xml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">X1 has a slope of - 0.05041. This means that </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">it </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>doesnt</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> really have an impact on whether y is </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>+1.</w:t></w:r>' +
  '</w:p>'
$para3.Range.InsertXML($pkgHeader + $p3xml + $pkgFooter)

# --- New paragraphs appended after "A(iii)" --------------------------------
$newParas =
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Boundry</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> is a linear plane and it </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>seoerates</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> the 1 and -1. It can be </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>witten</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> as </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>w1*x1 + w2*x2 + b = 0</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>w1, w2 = Weights.</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>x1, x</w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2  =</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> features</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>b = intercept</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">we use minus in the formula </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">to flip the sign of the coefficients. This is because the </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>descion</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>boundtry</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> is where the log- odds is equal to 0 .</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($pkgHeader + $newParas + $pkgFooter)
